$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Update C6 text on the main sheet: "ValidSearchHeader" -> "ValidSearchHeader1"
$ws1.Range("C6").Value = "ValidSearchHeader1"

# Move selection on sheet1 to C7
$ws1.Activate() | Out-Null
$ws1.Range("C7").Select() | Out-Null

# Move selection on sheet2 back to B6 (unchanged) - ensures sheet2 is not the active one anymore
$ws2.Range("B6").Select() | Out-Null

# Re-activate sheet1 so it becomes the workbook's active/visible tab
$ws1.Activate() | Out-Null
